$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 8 (shifts existing rows 8-11 down to 11-14,
# carrying their data & formatting along)
$ws.Rows("8:10").Insert()

# Fill in the 3 new rows (8, 9, 10) with the new weekly price data
$newRows = @(
    @{ Row=8;  D=44533; L='Primera'; M=350; N=24000; O=24000; P=24000; Q='$/caja 18 kilos'; R='Provincia de San Felipe de Aconcagua'; S=1333; T=18 },
    @{ Row=9;  D=44533; L='Segunda'; M=350; N=20000; O=20000; P=20000; Q='$/caja 18 kilos'; R='Provincia de San Felipe de Aconcagua'; S=1111; T=18 },
    @{ Row=10; D=44533; L='Tercera'; M=350; N=17000; O=17000; P=17000; Q='$/caja 18 kilos'; R='Provincia de San Felipe de Aconcagua'; S=944;  T=18 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 4
    $ws.Cells.Item($row, 2).Value = 'Feria Lagunitas de Puerto Montt'
    $ws.Cells.Item($row, 3).Value = 'Los Lagos'
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 10
    $ws.Cells.Item($row, 6).Value = 'Fruta'
    $ws.Cells.Item($row, 7).Value = 100103
    $ws.Cells.Item($row, 8).Value = 'Frutos de hueso (carozo)'
    $ws.Cells.Item($row, 9).Value = 100103003
    $ws.Cells.Item($row, 10).Value = 'Damasco'
    $ws.Cells.Item($row, 11).Value = 'Castle Brite'
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}

Write-Host "done"
